# Update cryptos list values (Price / Volume(1h), and a couple of row
# re-orderings for Uniswap/Chainlink and Aave/Cronos) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.907.15"
$ws.Range("E2").Value = "  +1.11%  "

$ws.Range("D3").Value = "2.906.97"
$ws.Range("E3").Value = "  +3.43%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "596.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  -1.18%  "

$ws.Range("E9").Value = "  +2.55%  "

$ws.Range("D10").Value = "2.905.00"
$ws.Range("E10").Value = "  +3.51%  "

$ws.Range("E11").Value = "  +15.00%  "

$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("D14").Value = "3.439.49"
$ws.Range("E14").Value = "  +3.30%  "

$ws.Range("D15").Value = "75.776.52"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000190"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.57%  "

$ws.Range("D18").Value = "2.916.61"
$ws.Range("E18").Value = "  +3.59%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.92%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.19%  "

$ws.Range("E23").Value = "  +4.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "3.059.10"
$ws.Range("E26").Value = "  +3.24%  "

$ws.Range("E27").Value = "  -0.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.49%  "

$ws.Range("E29").Value = "  +3.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.34%  "

$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "500.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.39%  "

$ws.Range("E34").Value = "  +0.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.51%  "

$ws.Range("E39").Value = "  -5.77%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "180.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.100"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +16.45%  "

$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.27%  "

$ws.Range("E45").Value = "  -2.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.31%  "

$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.570"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.654"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.24%  "
